$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 fixes ---
# B14: was inline text "2" -> becomes a real number 2
$ws.Range("B14").Value = 2

# C14: was inline text "nan" -> becomes an empty TEXT value (still present
# as a string cell, just with no characters). A bare "" assignment collapses
# the cell to fully blank (matches real Excel semantics), so instead enter a
# lone quote-prefix marker (forces a text cell with empty contents) and then
# restore the Normal style so no stray formatting diff is left behind.
$ws.Range("C14").Value = "'"
$ws.Range("C14").Style = "Normal"

# --- New row 15 ---
$ws.Range("A15").Value = "parisk"

# B15 must stay literal TEXT "5" (not be auto-converted to the number 5).
# Force text entry by temporarily switching the cell to Text number format,
# then restore the Normal style so no stray formatting diff is left behind.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "5"
$ws.Range("B15").Style = "Normal"

$ws.Range("C15").Value = "well carried，very through"
$ws.Range("D15").Value = "APC"
$ws.Range("E15").Value = "MET"
$ws.Range("F15").Value = "ea04c829-c996-4167-8585-03efb193cd41"
$ws.Range("G15").Value = "ByOExmWAb_annotated.xlsx"
$ws.Range("H15").Value = "The experiments were well carried through and very thorough."
